# ---------------------------------------------------------------------------
# "some exercises is done like data validation."
#
# Adds a second mini-exercise block below the existing salary calculator:
#   - a NEW ACTIVITY box (ADD/SUB/MUL/DIV) in B16:E20
#   - a GRADES box (NUM/NAME/SCORE/GENDER) in G15:J26, with per-gender
#     subtotal (average) rows built the way Excel's Data > Subtotal feature
#     would (outline groups + SUBTOTAL formulas), an AutoFilter over the
#     header/detail rows, and a whole-number (0-100) data validation rule
#     on the SCORE column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- a few helper colors (OLE BGR-ish "RGB" longs used by COM) -----------
$RED          = 255            # FF0000
$BLUE_FILL    = 15189940       # B4C7E7 (same family already used: theme4 tint .6)
$GREEN_FILL   = 9359785        # A9D18E (same family already used: theme9 tint .4)
$GRAY_FILL    = 13224393       # C9C9C9 (theme6 "Accent3" tint .4)
$GOLD_FILL    = 10086143       # FFE699 (theme7 "Accent4" tint .6)
$LIGHTGRAY_FONT = 15132391     # E7E6E6 (theme2 "Background2")

# ---- column widths ---------------------------------------------------------
$ws.Columns("H").ColumnWidth = 12.109375
$ws.Columns("J").ColumnWidth = 18.109375

# ===========================================================================
# Row 15 -- "GRADES" banner over the new table (G15:J15 merged)
# ===========================================================================
$rng = $ws.Range("G15:J15")
$rng.Font.Bold = $true
$rng.Interior.Color = $GREEN_FILL
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.Merge()
$ws.Range("G15").Value = "GRADES"

# ===========================================================================
# Row 16 -- headers
# ===========================================================================
# "NEW ACTIVITY" banner (B16:E16 merged)
$rng = $ws.Range("B16:E16")
$rng.Font.Bold = $true
$rng.Font.Color = $LIGHTGRAY_FONT
$rng.Interior.Color = $GRAY_FILL
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.Merge()
$ws.Range("B16").Value = "NEW ACTIVITY"

# NUM / NAME / SCORE / GENDER headers
$rng = $ws.Range("G16:J16")
$rng.Interior.Color = $GOLD_FILL
$rng.HorizontalAlignment = -4108   # xlCenter
$ws.Range("G16").Value = "NUM"
$ws.Range("H16").Value = "NAME"
$ws.Range("I16").Value = "SCORE"
$ws.Range("J16").Value = "GENDER"

# ===========================================================================
# Rows 17-20 -- ADD / SUB / MUL / DIV mini exercise
# ===========================================================================
$labels = @("ADD", "SUB", "MUL", "DIV")
for ($i = 0; $i -lt 4; $i++) {
    $r = 17 + $i
    $ws.Range("B$r").Value = $labels[$i]
}
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 7
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 3

$ws.Range("E17").Formula = "=C17+D17"
$ws.Range("E18").Formula = "=C18-D18"
$ws.Range("E19").Formula = "=C19*D19"
$ws.Range("E20").Formula = "=C20/D20"

$bRng = $ws.Range("B17:B20")
$bRng.Interior.Color = $RED
$bRng.Borders.LineStyle = 1
$bRng.Borders.Weight = 2
$bRng.HorizontalAlignment = -4108

$cdRng = $ws.Range("C17:D20")
$cdRng.Interior.Color = $BLUE_FILL
$cdRng.Borders.LineStyle = 1
$cdRng.Borders.Weight = 2
$cdRng.HorizontalAlignment = -4108

$eRng = $ws.Range("E17:E20")
$eRng.Interior.Color = $GREEN_FILL
$eRng.Borders.LineStyle = 1
$eRng.Borders.Weight = 2
$eRng.HorizontalAlignment = -4108

# ===========================================================================
# Rows 17-26 -- GRADES table (NUM / NAME / SCORE / GENDER) with Male/Female
# subtotals, like Data > Subtotal ("Average" of SCORE, grouped by GENDER)
# ===========================================================================
$names  = @("JOHN", "JONAS", "JONATHAN", "JACK", "JIM")
$scores = @(20, 100, 80, 50, 35)
for ($i = 0; $i -lt 5; $i++) {
    $r = 17 + $i
    $ws.Range("G$r").Value = $i + 1
    $ws.Range("H$r").Value = $names[$i]
    $ws.Range("I$r").Value = $scores[$i]
    $ws.Range("J$r").Value = "MALE"
}

$ws.Range("I22").Formula = "=SUBTOTAL(1,I17:I21)"
$ws.Range("J22").Value = "MALE Average"

$names2  = @("JANE", "JULIE")
$scores2 = @(62, 89)
for ($i = 0; $i -lt 2; $i++) {
    $r = 23 + $i
    $ws.Range("G$r").Value = $i + 6
    $ws.Range("H$r").Value = $names2[$i]
    $ws.Range("I$r").Value = $scores2[$i]
    $ws.Range("J$r").Value = "FEMALE"
}

$ws.Range("I25").Formula = "=SUBTOTAL(1,I23:I24)"
$ws.Range("J25").Value = "FEMALE Average"

$ws.Range("I26").Formula = "=SUBTOTAL(1,I17:I24)"
$ws.Range("J26").Value = "Grand Average"

# ---- formatting for the GRADES table --------------------------------------
# NUM column: rows 17-26 all share the same "light blue" fill, centered
$gRng = $ws.Range("G17:G26")
$gRng.Interior.Color = $BLUE_FILL
$gRng.HorizontalAlignment = -4108

# detail rows 17-21 (NAME/SCORE), centered, no fill
$hiRng = $ws.Range("H17:I21")
$hiRng.HorizontalAlignment = -4108

# MALE/FEMALE subtotal rows 22 (NAME/SCORE cols), centered, no fill
$hi22 = $ws.Range("H22:I22")
$hi22.HorizontalAlignment = -4108

# detail rows 23-26 (NAME/SCORE), centered, no fill
$hiRng2 = $ws.Range("H23:I26")
$hiRng2.HorizontalAlignment = -4108

# Average labels (J22, J25, J26) - bold
$avgRng = $ws.Range("J22")
$avgRng.Font.Bold = $true
$avgRng = $ws.Range("J25")
$avgRng.Font.Bold = $true
$avgRng = $ws.Range("J26")
$avgRng.Font.Bold = $true

# ===========================================================================
# AutoFilter + hidden _FilterDatabase defined name (like Data > Filter)
# ===========================================================================
$ws.Range("G16:J24").AutoFilter() | Out-Null
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sayfa1!`$G`$16:`$J`$24")
$fdb.Visible = $false

# ===========================================================================
# Data validation: SCORE column must be a whole number between 0 and 100
# ===========================================================================
$valRng = $ws.Range("I17:I21")
$valRng.Validation.Add(1, 1, 1, 0, 100) | Out-Null
$valRng.Validation.ErrorTitle = "Outside range error"
$valRng.Validation.ErrorMessage = "please enter between 0-100."

# ===========================================================================
# Outline / grouping: detail rows nested one level below the subtotal rows
# ===========================================================================
$ws.Rows("17:22").Group()
$ws.Rows("17:21").Group()
$ws.Rows("23:25").Group()
$ws.Rows("23:24").Group()

# ===========================================================================
# Move the stray helper formula from K27 down to K30 (room for the new rows)
# ===========================================================================
$ws.Range("K27").ClearContents()
$ws.Range("K30").Formula = "=G8"

# ===========================================================================
# View tidy-up: scroll down a bit and land the selection near the new table
# ===========================================================================
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("L22").Select()
